$wb = $excel.ActiveWorkbook

# Update the Approver sheet: replace the approver email address
$approverSheet = $wb.Worksheets.Item("Approver")
$approverSheet.Range("A2").Value = "lward@hl.com.test"

# Update the hyperlink target to match the new address
if ($approverSheet.Hyperlinks.Count -gt 0) {
    foreach ($hl in $approverSheet.Hyperlinks) {
        if ($hl.Address -eq "mailto:gksegura@hl.com.test") {
            $hl.Address = "mailto:lward@hl.com.test"
        }
    }
}

# Switch active sheet to Approver and move the selection to D11
$approverSheet.Activate() | Out-Null
$approverSheet.Range("D11").Select() | Out-Null
